# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Septiembre de 2020 a las 06:15"

# Row 5 - India
$ws.Range("B5").Value = 4562414
$ws.Range("C5").Value = 2689
$ws.Range("D5").Value = 3542663
$ws.Range("E5").Value = 943447

# Row 20 - Pakistan
$ws.Range("B20").Value = 300371
$ws.Range("C20").Value = 516
$ws.Range("D20").Value = 288206
$ws.Range("E20").Value = 5795
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 6370

# Row 33 - Kazajistan
$ws.Range("D33").Value = 100409
$ws.Range("E33").Value = 4618

# Row 50 - Honduras
$ws.Range("B50").Value = 65802
$ws.Range("C50").Value = 205
$ws.Range("D50").Value = 15577
$ws.Range("E50").Value = 48176
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 2049

# Row 126 - Tailandia
$ws.Range("B126").Value = 3461
$ws.Range("C126").Value = 7
$ws.Range("D126").Value = 3312
$ws.Range("E126").Value = 91

# Row 184 - Mongolia
$ws.Range("B184").Value = 311
$ws.Range("C184").Value = 1
$ws.Range("E184").Value = 13

# Row 187 - Butan
$ws.Range("B187").Value = 238
$ws.Range("C187").Value = 4
$ws.Range("E187").Value = 83
